# "update scripts wuth new tpm": NATMI re-run dropped the "Resolving-Mac"
# target-cluster rows (originally rows 6, 11 and 16) and refreshed every
# numeric column (E:T) for the remaining sending/target cluster pairs with
# values from the new TPM-based computation. Sending/Ligand/Receptor/Target
# (columns A:D) for the surviving rows are unchanged, so only structural row
# removal plus a value refresh are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Target cluster is "Resolving-Mac".
# Deleting from the bottom up means each index still refers to the
# original row because nothing below it has shifted yet.
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(6).Delete()

# Refreshed values for columns E (Ligand-expressing cells) through
# T (Edge total expression derived specificity) for the 12 remaining
# data rows (now rows 2-13).
$data = @(
    @(3, 1, 1.172733, 3.518199, 0.02094912533655296, 0.02094912533655295, 3, 1, 26.60444266666667, 79.813328, 0.5736225649467147, 0.5736225649467147, 31.199907861808, 280.799170756272, 0.01201689100894372, 0.01201689100894371),
    @(3, 1, 1.172733, 3.518199, 0.02094912533655296, 0.02094912533655295, 3, 1, 14.03147833333333, 42.094435, 0.3025349071358453, 0.3025349071358453, 16.455177680285, 148.096599122565, 0.006337841688271233, 0.006337841688271232),
    @(3, 1, 1.172733, 3.518199, 0.02094912533655296, 0.02094912533655295, 1, 0.3333333333333333, 0.2022123333333333, 0.606637, 0.004359931864156574, 0.004359931864156574, 0.237141076307, 2.134269686763, 0.00009133675908104706, 0.00009133675908104705),
    @(3, 1, 1.172733, 3.518199, 0.02094912533655296, 0.02094912533655295, 3, 1, 5.541567000000001, 16.624701, 0.1194825960532834, 0.1194825960532834, 6.498778492611001, 58.48900643349901, 0.002503055880256962, 0.002503055880256961),
    @(3, 1, 54.58029933333334, 163.740898, 0.974995614211059, 0.974995614211059, 3, 1, 26.60444266666667, 79.813328, 0.5736225649467147, 0.5736225649467147, 1452.078444343172, 13068.70599908854, 0.5592794850355451, 0.5592794850355451),
    @(3, 1, 54.58029933333334, 163.740898, 0.974995614211059, 0.974995614211059, 3, 1, 14.03147833333333, 42.094435, 0.3025349071358453, 0.3025349071358453, 765.8422875225145, 6892.580587702631, 0.2949702076031992, 0.2949702076031992),
    @(3, 1, 54.58029933333334, 163.740898, 0.974995614211059, 0.974995614211059, 1, 0.3333333333333333, 0.2022123333333333, 0.606637, 0.004359931864156574, 0.004359931864156574, 11.03680968222511, 99.33128714002601, 0.004250914445811707, 0.004250914445811707),
    @(3, 1, 54.58029933333334, 163.740898, 0.974995614211059, 0.974995614211059, 3, 1, 5.541567000000001, 16.624701, 0.1194825960532834, 0.1194825960532834, 302.460385635722, 2722.143470721498, 0.1164950071265029, 0.1164950071265029),
    @(1, 0.3333333333333333, 0.2270136666666667, 0.681041, 0.0040552604523881, 0.0040552604523881, 3, 1, 26.60444266666667, 79.813328, 0.5736225649467147, 0.5736225649467147, 6.039572079383111, 54.356148714448, 0.002326188902225836, 0.002326188902225836),
    @(1, 0.3333333333333333, 0.2270136666666667, 0.681041, 0.0040552604523881, 0.0040552604523881, 3, 1, 14.03147833333333, 42.094435, 0.3025349071358453, 0.3025349071358453, 3.185337345203889, 28.668036106835, 0.0012268578443749, 0.0012268578443749),
    @(1, 0.3333333333333333, 0.2270136666666667, 0.681041, 0.0040552604523881, 0.0040552604523881, 1, 0.3333333333333333, 0.2022123333333333, 0.606637, 0.004359931864156574, 0.004359931864156574, 0.04590496323522222, 0.413144669117, 0.00001768065926382088, 0.00001768065926382088),
    @(1, 0.3333333333333333, 0.2270136666666667, 0.681041, 0.0040552604523881, 0.0040552604523881, 3, 1, 5.541567000000001, 16.624701, 0.1194825960532834, 0.1194825960532834, 1.258011443749, 11.322102993741, 0.0004845330465235426, 0.0004845330465235427)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = 2 + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, 5 + $j).Value = $row[$j]
    }
}
